# ---- 1) Insert new "2022-Q4" sheet before "2022-Q3" and populate with data ----
$wb = $excel.ActiveWorkbook
$targetSheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($targetSheet)
$newSheet.Name = "2022-Q4"

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'166019"
$newSheet.Cells.Item(2,3).Value = "中欧价值智选回报混合A"
$newSheet.Cells.Item(2,4).Value = "'71.44"
$newSheet.Cells.Item(2,5).Value = "'94.96"
$newSheet.Cells.Item(2,6).Value = "'3.35"
$newSheet.Cells.Item(2,7).Value = "'2.3932"
$newSheet.Cells.Item(2,8).Value = 10
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'013220"
$newSheet.Cells.Item(3,3).Value = "中欧新兴价值一年持有混合A"
$newSheet.Cells.Item(3,4).Value = "'36.58"
$newSheet.Cells.Item(3,5).Value = "'93.35"
$newSheet.Cells.Item(3,6).Value = "'3.36"
$newSheet.Cells.Item(3,7).Value = "'1.2291"
$newSheet.Cells.Item(3,8).Value = 9
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'010363"
$newSheet.Cells.Item(4,3).Value = "信澳匠心臻选两年持有期混合"
$newSheet.Cells.Item(4,4).Value = "'37.44"
$newSheet.Cells.Item(4,5).Value = "'92.71"
$newSheet.Cells.Item(4,6).Value = "'3.10"
$newSheet.Cells.Item(4,7).Value = "'1.1606"
$newSheet.Cells.Item(4,8).Value = 6
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "'004235"
$newSheet.Cells.Item(5,3).Value = "中欧价值智选回报混合C"
$newSheet.Cells.Item(5,4).Value = "'29.51"
$newSheet.Cells.Item(5,5).Value = "'94.96"
$newSheet.Cells.Item(5,6).Value = "'3.35"
$newSheet.Cells.Item(5,7).Value = "'0.9886"
$newSheet.Cells.Item(5,8).Value = 10
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "'013221"
$newSheet.Cells.Item(6,3).Value = "中欧新兴价值一年持有混合C"
$newSheet.Cells.Item(6,4).Value = "'16.10"
$newSheet.Cells.Item(6,5).Value = "'93.35"
$newSheet.Cells.Item(6,6).Value = "'3.36"
$newSheet.Cells.Item(6,7).Value = "'0.5410"
$newSheet.Cells.Item(6,8).Value = 9
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "'014404"
$newSheet.Cells.Item(7,3).Value = "中欧多元价值三年持有混合A"
$newSheet.Cells.Item(7,4).Value = "'14.47"
$newSheet.Cells.Item(7,5).Value = "'91.65"
$newSheet.Cells.Item(7,6).Value = "'3.37"
$newSheet.Cells.Item(7,7).Value = "'0.4876"
$newSheet.Cells.Item(7,8).Value = 9
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "'009011"
$newSheet.Cells.Item(8,3).Value = "华夏睿阳一年持有期混合"
$newSheet.Cells.Item(8,4).Value = "'14.05"
$newSheet.Cells.Item(8,5).Value = "'82.02"
$newSheet.Cells.Item(8,6).Value = "'2.67"
$newSheet.Cells.Item(8,7).Value = "'0.3751"
$newSheet.Cells.Item(8,8).Value = 4
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "'001887"
$newSheet.Cells.Item(9,3).Value = "中欧价值智选回报混合E"
$newSheet.Cells.Item(9,4).Value = "'9.15"
$newSheet.Cells.Item(9,5).Value = "'94.96"
$newSheet.Cells.Item(9,6).Value = "'3.35"
$newSheet.Cells.Item(9,7).Value = "'0.3065"
$newSheet.Cells.Item(9,8).Value = 10
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "'610004"
$newSheet.Cells.Item(10,3).Value = "信澳中小盘混合"
$newSheet.Cells.Item(10,4).Value = "'4.60"
$newSheet.Cells.Item(10,5).Value = "'91.92"
$newSheet.Cells.Item(10,6).Value = "'5.61"
$newSheet.Cells.Item(10,7).Value = "'0.2581"
$newSheet.Cells.Item(10,8).Value = 6
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "'001105"
$newSheet.Cells.Item(11,3).Value = "信澳转型创新股票"
$newSheet.Cells.Item(11,4).Value = "'3.06"
$newSheet.Cells.Item(11,5).Value = "'93.59"
$newSheet.Cells.Item(11,6).Value = "'3.50"
$newSheet.Cells.Item(11,7).Value = "'0.1071"
$newSheet.Cells.Item(11,8).Value = 7
$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = "'004119"
$newSheet.Cells.Item(12,3).Value = "广发创新驱动灵活配置混合"
$newSheet.Cells.Item(12,4).Value = "'1.78"
$newSheet.Cells.Item(12,5).Value = "'74.78"
$newSheet.Cells.Item(12,6).Value = "'4.97"
$newSheet.Cells.Item(12,7).Value = "'0.0885"
$newSheet.Cells.Item(12,8).Value = 7
$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = "'013495"
$newSheet.Cells.Item(13,3).Value = "信澳产业优选一年持有混合A"
$newSheet.Cells.Item(13,4).Value = "'1.92"
$newSheet.Cells.Item(13,5).Value = "'79.53"
$newSheet.Cells.Item(13,6).Value = "'4.52"
$newSheet.Cells.Item(13,7).Value = "'0.0868"
$newSheet.Cells.Item(13,8).Value = 8
$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = "'007146"
$newSheet.Cells.Item(14,3).Value = "鹏华研究智选混合"
$newSheet.Cells.Item(14,4).Value = "'3.79"
$newSheet.Cells.Item(14,5).Value = "'79.83"
$newSheet.Cells.Item(14,6).Value = "'1.72"
$newSheet.Cells.Item(14,7).Value = "'0.0652"
$newSheet.Cells.Item(14,8).Value = 7
$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,2).Value = "'014135"
$newSheet.Cells.Item(15,3).Value = "中欧金安量化混合A"
$newSheet.Cells.Item(15,4).Value = "'7.09"
$newSheet.Cells.Item(15,5).Value = "'90.01"
$newSheet.Cells.Item(15,6).Value = "'0.88"
$newSheet.Cells.Item(15,7).Value = "'0.0624"
$newSheet.Cells.Item(15,8).Value = 1
$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,2).Value = "'005683"
$newSheet.Cells.Item(16,3).Value = "国寿安保华兴灵活配置混合"
$newSheet.Cells.Item(16,4).Value = "'2.11"
$newSheet.Cells.Item(16,5).Value = "'92.97"
$newSheet.Cells.Item(16,6).Value = "'2.80"
$newSheet.Cells.Item(16,7).Value = "'0.0591"
$newSheet.Cells.Item(16,8).Value = 3
$newSheet.Cells.Item(17,1).Value = 15
$newSheet.Cells.Item(17,2).Value = "'880007"
$newSheet.Cells.Item(17,3).Value = "招商资管智远成长灵活配置混合A"
$newSheet.Cells.Item(17,4).Value = "'1.49"
$newSheet.Cells.Item(17,5).Value = "'90.93"
$newSheet.Cells.Item(17,6).Value = "'3.64"
$newSheet.Cells.Item(17,7).Value = "'0.0542"
$newSheet.Cells.Item(17,8).Value = 7
$newSheet.Cells.Item(18,1).Value = 16
$newSheet.Cells.Item(18,2).Value = "'014405"
$newSheet.Cells.Item(18,3).Value = "中欧多元价值三年持有混合C"
$newSheet.Cells.Item(18,4).Value = "'1.36"
$newSheet.Cells.Item(18,5).Value = "'91.65"
$newSheet.Cells.Item(18,6).Value = "'3.37"
$newSheet.Cells.Item(18,7).Value = "'0.0458"
$newSheet.Cells.Item(18,8).Value = 9
$newSheet.Cells.Item(19,1).Value = 17
$newSheet.Cells.Item(19,2).Value = "'003131"
$newSheet.Cells.Item(19,3).Value = "国寿安保强国智造灵活配置混合"
$newSheet.Cells.Item(19,4).Value = "'1.38"
$newSheet.Cells.Item(19,5).Value = "'92.54"
$newSheet.Cells.Item(19,6).Value = "'3.27"
$newSheet.Cells.Item(19,7).Value = "'0.0451"
$newSheet.Cells.Item(19,8).Value = 3
$newSheet.Cells.Item(20,1).Value = 18
$newSheet.Cells.Item(20,2).Value = "'004332"
$newSheet.Cells.Item(20,3).Value = "恒生前海沪港深新兴产业精选混合"
$newSheet.Cells.Item(20,4).Value = "'0.49"
$newSheet.Cells.Item(20,5).Value = "'75.81"
$newSheet.Cells.Item(20,6).Value = "'6.95"
$newSheet.Cells.Item(20,7).Value = "'0.0341"
$newSheet.Cells.Item(20,8).Value = 2
$newSheet.Cells.Item(21,1).Value = 19
$newSheet.Cells.Item(21,2).Value = "'001990"
$newSheet.Cells.Item(21,3).Value = "中欧数据挖掘多因子灵活配置混合A"
$newSheet.Cells.Item(21,4).Value = "'3.26"
$newSheet.Cells.Item(21,5).Value = "'90.14"
$newSheet.Cells.Item(21,6).Value = "'0.89"
$newSheet.Cells.Item(21,7).Value = "'0.0290"
$newSheet.Cells.Item(21,8).Value = 1
$newSheet.Cells.Item(22,1).Value = 20
$newSheet.Cells.Item(22,2).Value = "'011735"
$newSheet.Cells.Item(22,3).Value = "国寿安保裕丰混合C"
$newSheet.Cells.Item(22,4).Value = "'2.64"
$newSheet.Cells.Item(22,5).Value = "'27.90"
$newSheet.Cells.Item(22,6).Value = "'0.85"
$newSheet.Cells.Item(22,7).Value = "'0.0224"
$newSheet.Cells.Item(22,8).Value = 4
$newSheet.Cells.Item(23,1).Value = 21
$newSheet.Cells.Item(23,2).Value = "'881007"
$newSheet.Cells.Item(23,3).Value = "招商资管智远成长灵活配置混合C"
$newSheet.Cells.Item(23,4).Value = "'0.60"
$newSheet.Cells.Item(23,5).Value = "'90.93"
$newSheet.Cells.Item(23,6).Value = "'3.64"
$newSheet.Cells.Item(23,7).Value = "'0.0218"
$newSheet.Cells.Item(23,8).Value = 7
$newSheet.Cells.Item(24,1).Value = 22
$newSheet.Cells.Item(24,2).Value = "'010206"
$newSheet.Cells.Item(24,3).Value = "国寿安保裕安混合C"
$newSheet.Cells.Item(24,4).Value = "'2.00"
$newSheet.Cells.Item(24,5).Value = "'33.78"
$newSheet.Cells.Item(24,6).Value = "'1.03"
$newSheet.Cells.Item(24,7).Value = "'0.0206"
$newSheet.Cells.Item(24,8).Value = 3
$newSheet.Cells.Item(25,1).Value = 23
$newSheet.Cells.Item(25,2).Value = "'970113"
$newSheet.Cells.Item(25,3).Value = "兴证资管金麒麟兴睿优选一年持有期混合B"
$newSheet.Cells.Item(25,4).Value = "'0.67"
$newSheet.Cells.Item(25,5).Value = "'84.89"
$newSheet.Cells.Item(25,6).Value = "'3.02"
$newSheet.Cells.Item(25,7).Value = "'0.0202"
$newSheet.Cells.Item(25,8).Value = 9
$newSheet.Cells.Item(26,1).Value = 24
$newSheet.Cells.Item(26,2).Value = "'010205"
$newSheet.Cells.Item(26,3).Value = "国寿安保裕安混合A"
$newSheet.Cells.Item(26,4).Value = "'1.89"
$newSheet.Cells.Item(26,5).Value = "'33.78"
$newSheet.Cells.Item(26,6).Value = "'1.03"
$newSheet.Cells.Item(26,7).Value = "'0.0195"
$newSheet.Cells.Item(26,8).Value = 3
$newSheet.Cells.Item(27,1).Value = 25
$newSheet.Cells.Item(27,2).Value = "'004234"
$newSheet.Cells.Item(27,3).Value = "中欧数据挖掘多因子灵活配置混合C"
$newSheet.Cells.Item(27,4).Value = "'1.93"
$newSheet.Cells.Item(27,5).Value = "'90.14"
$newSheet.Cells.Item(27,6).Value = "'0.89"
$newSheet.Cells.Item(27,7).Value = "'0.0172"
$newSheet.Cells.Item(27,8).Value = 1
$newSheet.Cells.Item(28,1).Value = 26
$newSheet.Cells.Item(28,2).Value = "'001420"
$newSheet.Cells.Item(28,3).Value = "南方大数据300指数A"
$newSheet.Cells.Item(28,4).Value = "'1.71"
$newSheet.Cells.Item(28,5).Value = "'93.56"
$newSheet.Cells.Item(28,6).Value = "'1.00"
$newSheet.Cells.Item(28,7).Value = "'0.0171"
$newSheet.Cells.Item(28,8).Value = 2
$newSheet.Cells.Item(29,1).Value = 27
$newSheet.Cells.Item(29,2).Value = "'016370"
$newSheet.Cells.Item(29,3).Value = "信澳业绩驱动混合A"
$newSheet.Cells.Item(29,4).Value = "'0.77"
$newSheet.Cells.Item(29,5).Value = "'30.31"
$newSheet.Cells.Item(29,6).Value = "'1.75"
$newSheet.Cells.Item(29,7).Value = "'0.0135"
$newSheet.Cells.Item(29,8).Value = 7
$newSheet.Cells.Item(30,1).Value = 28
$newSheet.Cells.Item(30,2).Value = "'015608"
$newSheet.Cells.Item(30,3).Value = "信澳转型创新股票C"
$newSheet.Cells.Item(30,4).Value = "'0.33"
$newSheet.Cells.Item(30,5).Value = "'93.59"
$newSheet.Cells.Item(30,6).Value = "'3.50"
$newSheet.Cells.Item(30,7).Value = "'0.0116"
$newSheet.Cells.Item(30,8).Value = 7
$newSheet.Cells.Item(31,1).Value = 29
$newSheet.Cells.Item(31,2).Value = "'011734"
$newSheet.Cells.Item(31,3).Value = "国寿安保裕丰混合A"
$newSheet.Cells.Item(31,4).Value = "'1.22"
$newSheet.Cells.Item(31,5).Value = "'27.90"
$newSheet.Cells.Item(31,6).Value = "'0.85"
$newSheet.Cells.Item(31,7).Value = "'0.0104"
$newSheet.Cells.Item(31,8).Value = 4
$newSheet.Cells.Item(32,1).Value = 30
$newSheet.Cells.Item(32,2).Value = "'013383"
$newSheet.Cells.Item(32,3).Value = "恒生前海高端制造混合A"
$newSheet.Cells.Item(32,4).Value = "'0.11"
$newSheet.Cells.Item(32,5).Value = "'84.98"
$newSheet.Cells.Item(32,6).Value = "'9.15"
$newSheet.Cells.Item(32,7).Value = "'0.0101"
$newSheet.Cells.Item(32,8).Value = 1
$newSheet.Cells.Item(33,1).Value = 31
$newSheet.Cells.Item(33,2).Value = "'014136"
$newSheet.Cells.Item(33,3).Value = "中欧金安量化混合C"
$newSheet.Cells.Item(33,4).Value = "'1.07"
$newSheet.Cells.Item(33,5).Value = "'90.01"
$newSheet.Cells.Item(33,6).Value = "'0.88"
$newSheet.Cells.Item(33,7).Value = "'0.0094"
$newSheet.Cells.Item(33,8).Value = 1
$newSheet.Cells.Item(34,1).Value = 32
$newSheet.Cells.Item(34,2).Value = "'013496"
$newSheet.Cells.Item(34,3).Value = "信澳产业优选一年持有混合C"
$newSheet.Cells.Item(34,4).Value = "'0.18"
$newSheet.Cells.Item(34,5).Value = "'79.53"
$newSheet.Cells.Item(34,6).Value = "'4.52"
$newSheet.Cells.Item(34,7).Value = "'0.0081"
$newSheet.Cells.Item(34,8).Value = 8
$newSheet.Cells.Item(35,1).Value = 33
$newSheet.Cells.Item(35,2).Value = "'016371"
$newSheet.Cells.Item(35,3).Value = "信澳业绩驱动混合C"
$newSheet.Cells.Item(35,4).Value = "'0.25"
$newSheet.Cells.Item(35,5).Value = "'30.31"
$newSheet.Cells.Item(35,6).Value = "'1.75"
$newSheet.Cells.Item(35,7).Value = "'0.0044"
$newSheet.Cells.Item(35,8).Value = 7
$newSheet.Cells.Item(36,1).Value = 34
$newSheet.Cells.Item(36,2).Value = "'011771"
$newSheet.Cells.Item(36,3).Value = "国寿安保稳隆混合A"
$newSheet.Cells.Item(36,4).Value = "'0.50"
$newSheet.Cells.Item(36,5).Value = "'32.85"
$newSheet.Cells.Item(36,6).Value = "'0.82"
$newSheet.Cells.Item(36,7).Value = "'0.0041"
$newSheet.Cells.Item(36,8).Value = 6
$newSheet.Cells.Item(37,1).Value = 35
$newSheet.Cells.Item(37,2).Value = "'013384"
$newSheet.Cells.Item(37,3).Value = "恒生前海高端制造混合C"
$newSheet.Cells.Item(37,4).Value = "'0.04"
$newSheet.Cells.Item(37,5).Value = "'84.98"
$newSheet.Cells.Item(37,6).Value = "'9.15"
$newSheet.Cells.Item(37,7).Value = "'0.0037"
$newSheet.Cells.Item(37,8).Value = 1
$newSheet.Cells.Item(38,1).Value = 36
$newSheet.Cells.Item(38,2).Value = "'001426"
$newSheet.Cells.Item(38,3).Value = "南方大数据300指数C"
$newSheet.Cells.Item(38,4).Value = "'0.32"
$newSheet.Cells.Item(38,5).Value = "'93.56"
$newSheet.Cells.Item(38,6).Value = "'1.00"
$newSheet.Cells.Item(38,7).Value = "'0.0032"
$newSheet.Cells.Item(38,8).Value = 2
$newSheet.Cells.Item(39,1).Value = 37
$newSheet.Cells.Item(39,2).Value = "'001932"
$newSheet.Cells.Item(39,3).Value = "国寿安保灵活优选混合"
$newSheet.Cells.Item(39,4).Value = "'0.11"
$newSheet.Cells.Item(39,5).Value = "'39.50"
$newSheet.Cells.Item(39,6).Value = "'1.95"
$newSheet.Cells.Item(39,7).Value = "'0.0021"
$newSheet.Cells.Item(39,8).Value = 1
$newSheet.Cells.Item(40,1).Value = 38
$newSheet.Cells.Item(40,2).Value = "'008533"
$newSheet.Cells.Item(40,3).Value = "惠升惠兴混合A"
$newSheet.Cells.Item(40,4).Value = "'0.03"
$newSheet.Cells.Item(40,5).Value = "'23.29"
$newSheet.Cells.Item(40,6).Value = "'1.29"
$newSheet.Cells.Item(40,7).Value = "'0.0004"
$newSheet.Cells.Item(40,8).Value = 9
$newSheet.Cells.Item(41,1).Value = 39
$newSheet.Cells.Item(41,2).Value = "'970112"
$newSheet.Cells.Item(41,3).Value = "兴证资管金麒麟兴睿优选一年持有期混合A"
$newSheet.Cells.Item(41,4).Value = "'0.00"
$newSheet.Cells.Item(41,5).Value = "'84.89"
$newSheet.Cells.Item(41,6).Value = "'3.02"
$newSheet.Cells.Item(41,7).Value = 0
$newSheet.Cells.Item(41,8).Value = 9
$newSheet.Cells.Item(42,1).Value = 40
$newSheet.Cells.Item(42,2).Value = "'970114"
$newSheet.Cells.Item(42,3).Value = "兴证资管金麒麟兴睿优选一年持有期混合C"
$newSheet.Cells.Item(42,4).Value = "'0.00"
$newSheet.Cells.Item(42,5).Value = "'84.89"
$newSheet.Cells.Item(42,6).Value = "'3.02"
$newSheet.Cells.Item(42,7).Value = 0
$newSheet.Cells.Item(42,8).Value = 9
$newSheet.Cells.Item(43,1).Value = 41
$newSheet.Cells.Item(43,2).Value = "'011772"
$newSheet.Cells.Item(43,3).Value = "国寿安保稳隆混合C"
$newSheet.Cells.Item(43,4).Value = "'0.00"
$newSheet.Cells.Item(43,5).Value = "'32.85"
$newSheet.Cells.Item(43,6).Value = "'0.82"
$newSheet.Cells.Item(43,7).Value = 0
$newSheet.Cells.Item(43,8).Value = 6
$newSheet.Cells.Item(44,1).Value = 42
$newSheet.Cells.Item(44,2).Value = "'008534"
$newSheet.Cells.Item(44,3).Value = "惠升惠兴混合C"
$newSheet.Cells.Item(44,4).Value = "'0.00"
$newSheet.Cells.Item(44,5).Value = "'23.29"
$newSheet.Cells.Item(44,6).Value = "'1.29"
$newSheet.Cells.Item(44,7).Value = 0
$newSheet.Cells.Item(44,8).Value = 9

# Fix formatting: clear the quote-prefix formatting picked up from text entry,
# then re-apply the workbook's header/index style (matches the other quarter sheets).
$newSheet.Range("A1:H44").ClearFormats()
$styleSource = $wb.Worksheets.Item("2022-Q3").Range("B1")
$styleSource.Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("A2:A44").PasteSpecial(-4122)
$excel.CutCopyMode = $false
# ---- 2) Update the "总计" (summary) sheet: insert a new top data row for 2022-Q4 ----
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Re-apply the index-column style (border+bold) that Insert() doesn't carry
# into the freshly-created row, and strip the stray formatting it does add
# to the other cells so the new row matches its siblings.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()
$excel.CutCopyMode = $false

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 43
$summary.Cells.Item(2,4).Value = 8.640000000000001

# Renumber the index column for the rows that shifted down one position.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r,1).Value = $r - 2
}
